$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 29: update title (D29), link (E29) unchanged
$ws.Range("D29").Value = "프로메디우스"

# Row 37: update title (D37) and link (E37)
$ws.Range("D37").Value = "[Rehearsal] 석사학위논문심사 - 김명섭"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1857&mod=document&pageid=1"

# Row 41: update title (D41) and link (E41)
$ws.Range("D41").Value = "양자 컴퓨팅과 미래"
$ws.Range("E41").Value = "http://cloudinsight.net/cloud/%ec%96%91%ec%9e%90-%ec%bb%b4%ed%93%a8%ed%8c%85%ea%b3%bc-%eb%af%b8%eb%9e%98/"

# Row 50: update title (D50) and link (E50)
$ws.Range("D50").Value = "GPU 장비의 중요성 그리고 시장의 반응"
$ws.Range("E50").Value = "http://incredible.egloos.com/7530150"
